{"js": "// Update the worksheet date and the 25 three-digit-by-one-digit\n// multiplication prompts to the new values from the commit.\nconst replacements = [\n  [\"2025-12-12 Friday\", \"2025-12-13 Saturday\"],\n  [\"924\u00d78=\", \"379\u00d73=\"],\n  [\"157\u00d73=\", \"792\u00d74=\"],\n  [\"136\u00d76=\", \"732\u00d76=\"],\n  [\"708\u00d77=\", \"900\u00d75=\"],\n  [\"589\u00d73=\", \"349\u00d73=\"],\n  [\"722\u00d73=\", \"585\u00d74=\"],\n  [\"776\u00d77=\", \"618\u00d79=\"],\n  [\"767\u00d72=\", \"961\u00d73=\"],\n  [\"308\u00d77=\", \"390\u00d77=\"],\n  [\"372\u00d72=\", \"905\u00d76=\"],\n  [\"161\u00d78=\", \"406\u00d78=\"],\n  [\"803\u00d76=\", \"238\u00d78=\"],\n  [\"950\u00d74=\", \"116\u00d78=\"],\n  [\"119\u00d75=\", \"992\u00d77=\"],\n  [\"484\u00d75=\", \"292\u00d75=\"],\n  [\"949\u00d74=\", \"929\u00d76=\"],\n  [\"211\u00d74=\", \"518\u00d77=\"],\n  [\"551\u00d78=\", \"295\u00d78=\"],\n  [\"791\u00d79=\", \"927\u00d74=\"],\n  [\"318\u00d75=\", \"554\u00d78=\"],\n  [\"168\u00d77=\", \"225\u00d77=\"],\n  [\"722\u00d72=\", \"755\u00d76=\"],\n  [\"763\u00d74=\", \"211\u00d77=\"],\n  [\"779\u00d79=\", \"817\u00d76=\"],\n  [\"390\u00d79=\", \"212\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 three-digit-by-one-digit\n# multiplication prompts to the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-12 Friday\", \"2025-12-13 Saturday\"),\n    @(\"924\u00d78=\", \"379\u00d73=\"),\n    @(\"157\u00d73=\", \"792\u00d74=\"),\n    @(\"136\u00d76=\", \"732\u00d76=\"),\n    @(\"708\u00d77=\", \"900\u00d75=\"),\n    @(\"589\u00d73=\", \"349\u00d73=\"),\n    @(\"722\u00d73=\", \"585\u00d74=\"),\n    @(\"776\u00d77=\", \"618\u00d79=\"),\n    @(\"767\u00d72=\", \"961\u00d73=\"),\n    @(\"308\u00d77=\", \"390\u00d77=\"),\n    @(\"372\u00d72=\", \"905\u00d76=\"),\n    @(\"161\u00d78=\", \"406\u00d78=\"),\n    @(\"803\u00d76=\", \"238\u00d78=\"),\n    @(\"950\u00d74=\", \"116\u00d78=\"),\n    @(\"119\u00d75=\", \"992\u00d77=\"),\n    @(\"484\u00d75=\", \"292\u00d75=\"),\n    @(\"949\u00d74=\", \"929\u00d76=\"),\n    @(\"211\u00d74=\", \"518\u00d77=\"),\n    @(\"551\u00d78=\", \"295\u00d78=\"),\n    @(\"791\u00d79=\", \"927\u00d74=\"),\n    @(\"318\u00d75=\", \"554\u00d78=\"),\n    @(\"168\u00d77=\", \"225\u00d77=\"),\n    @(\"722\u00d72=\", \"755\u00d76=\"),\n    @(\"763\u00d74=\", \"211\u00d77=\"),\n    @(\"779\u00d79=\", \"817\u00d76=\"),\n    @(\"390\u00d79=\", \"212\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
